$d = $word.ActiveDocument

# 1. Remove "Jeffrey Saylor" from the author byline on the title page.
$d.Content.Find.Execute(", Ollie Peel, Jeffrey Saylor", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", Ollie Peel", 2)

# 2. Remove the "Jeffrey Saylor: " contribution paragraph in its entirety
#    (its own paragraph mark and the mark ending the preceding "Ollie
#    Peel: " paragraph), so "Ollie Peel: " merges directly with the
#    paragraph that follows (which only holds the page-break run).
$d.Content.Find.Execute("`rJeffrey Saylor: `r", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
